$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns A-C become one uniform width, column D gets a new (narrower) width ---
$ws.Range("A1:C1").ColumnWidth = 31.57
$ws.Range("D1").ColumnWidth = 9.71

# --- Row heights: header rows get shorter, row 5 gets an explicit custom height ---
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 16.5

# --- New column T: 2023 data, one new value per existing data row, formatted like column S ---
$newValues = @{
    4  = 2023
    5  = 4.8
    6  = 5.7
    7  = 1.9
    8  = 8.9
    9  = 11.9
    10 = 2.5
    11 = 0.7
    12 = 12.7
    13 = 1.1
    14 = 2.2
}

foreach ($row in 4..14) {
    $src = $ws.Range("S$row")
    $dst = $ws.Range("T$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $newValues[$row]
}

$excel.CutCopyMode = 0

# --- Reset the active selection back to the default top-left cell ---
$ws.Range("A1").Select()
